# Final code with comment - tut 6
# Mark attendance counters (Real / Absent) for each date row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count + Real
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count + Real
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent
$ws.Range("H6").Value = 1

# Row 7: Absent
$ws.Range("H7").Value = 1

# Row 8: Absent
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count + Real
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Absent
$ws.Range("H10").Value = 1

# Row 11: Absent
$ws.Range("H11").Value = 1

# Row 12: Absent
$ws.Range("H12").Value = 1

# Row 13: Absent
$ws.Range("H13").Value = 1

# Row 14: Absent
$ws.Range("H14").Value = 1

# Row 15: Absent
$ws.Range("H15").Value = 1

# Row 16: Absent
$ws.Range("H16").Value = 1

# Row 17: Absent
$ws.Range("H17").Value = 1

# Row 18: Absent
$ws.Range("H18").Value = 1
